$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of raw data (A,B,C,D) to append starting at row 314
$data = @(
    @(4, 3, 5, 17),
    @(6, 15, 5, 5),
    @(2, 6, 3, 14),
    @(2, 14, 4, 6),
    @(8, 17, 5, 3),
    @(4, 17, 5, 3),
    @(4, 8, 5, 12),
    @(5, 8, 4, 12),
    @(1, 15, 2, 5),
    @(4, 13, 3, 7),
    @(8, 17, 5, 3),
    @(5, 12, 7, 8),
    @(5, 4, 6, 16),
    @(5, 4, 6, 16),
    @(3, 17, 2, 3),
    @(1, 16, 4, 4),
    @(5, 16, 4, 4)
)

$startRow = 314
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Assigning a formula to a whole multi-cell range in one shot makes the
# engine emit a proper shared-formula group (anchor cell carries the
# ref + formula text, the rest just carry the shared si). Rows 314-321
# continue the "B+D" pattern used by the block above them; rows 322-330
# are a second, freshly-entered block, so they land in their own group,
# just like Excel does when a new shared-formula range is started.
$ws.Range("E314:E321").Formula = "=B314+D314"
$ws.Range("E322:E330").Formula = "=B322+D322"

$ws.Range("A331").Select()
